$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 174, shifting existing rows 174-212 down to 175-213.
$ws.Rows("174:174").Insert()

# Populate the newly inserted row 174 with the new weekly price record.
$ws.Cells.Item(174, 1).Value = 10
$ws.Cells.Item(174, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(174, 3).Value = "La Araucanía"
$ws.Cells.Item(174, 4).Value = 44889
$ws.Cells.Item(174, 5).Value = 9
$ws.Cells.Item(174, 6).Value = 100112012
$ws.Cells.Item(174, 7).Value = "Espinaca"
$ws.Cells.Item(174, 8).Value = "Sin especificar"
$ws.Cells.Item(174, 9).Value = "Primera"
$ws.Cells.Item(174, 10).Value = 75
$ws.Cells.Item(174, 11).Value = 9000
$ws.Cells.Item(174, 12).Value = 9000
$ws.Cells.Item(174, 13).Value = 9000
$ws.Cells.Item(174, 14).Value = "$/docena de atados"
$ws.Cells.Item(174, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(174, 16).Value = 3000
$ws.Cells.Item(174, 17).Value = 3
$ws.Cells.Item(174, 18).Value = "Hortaliza"
